# Add data for 2021-11-06
# - Rename the worksheet title (and workbook sheet tab name) from
#   "Through 2021-10-28" to "Through 2021-10-29"
# - Update the "October (through 10-28)" label to "October (through 10-29)"
# - Update August 2021 total (H9) from 159 to 160
# - Update the October row (row 11) with the new day's data
# - Recompute the Total row (row 12) to reflect the updated values

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet tab to match the new "through" date.
$ws.Name = "Through 2021-10-29"

# Update the row label for October.
$ws.Range("A11").Value = "October (through 10-29)"

# August 2021 (row 9) value bump.
$ws.Range("H9").Value = 160

# October (row 11) updated values per column/year.
$ws.Range("B11").Value = 29
$ws.Range("C11").Value = 53
$ws.Range("D11").Value = 78
$ws.Range("E11").Value = 60
$ws.Range("F11").Value = 57
$ws.Range("G11").Value = 139
$ws.Range("H11").Value = 181

# Total row (row 12) recomputed sums.
$ws.Range("B12").Value = 255
$ws.Range("C12").Value = 482
$ws.Range("D12").Value = 705
$ws.Range("E12").Value = 608
$ws.Range("F12").Value = 479
$ws.Range("G12").Value = 1040
$ws.Range("H12").Value = 1429
